# Auto-generated edit script: apply precise cell-level changes
# derived from the target OOXML diff (row reshuffle + refresh-date bump).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 (source row 2) ---
$ws.Range("C2").Value = 46063
# --- row 3 (source row 3) ---
$ws.Range("C3").Value = 46063
# --- row 4 (source row 5) ---
$ws.Range("A4").Value = 'A 24802-2025'
$ws.Range("C4").Value = 46063
$ws.Range("G4").Value = 2.4
$ws.Range("R4").Value = 'Talltita'
$ws.Range("S4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/artfynd/A 24802-2025 artfynd.xlsx", "A 24802-2025")'
$ws.Range("T4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/kartor/A 24802-2025 karta.png", "A 24802-2025")'
$ws.Range("V4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/klagomål/A 24802-2025 FSC-klagomål.docx", "A 24802-2025")'
$ws.Range("W4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/klagomålsmail/A 24802-2025 FSC-klagomål mail.docx", "A 24802-2025")'
$ws.Range("X4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/tillsyn/A 24802-2025 tillsynsbegäran.docx", "A 24802-2025")'
$ws.Range("Y4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/tillsynsmail/A 24802-2025 tillsynsbegäran mail.docx", "A 24802-2025")'
$ws.Range("Z4").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/fåglar/A 24802-2025 prioriterade fågelarter.docx", "A 24802-2025")'
# --- row 5 (source row 6) ---
$ws.Range("A5").Value = 'A 1510-2024'
$ws.Range("B5").Value = 45306
$ws.Range("C5").Value = 46063
$ws.Range("G5").Value = 3.7
$ws.Range("J5").Value = 0
$ws.Range("N5").Value = 1
$ws.Range("R5").Value = 'Svart stork'
$ws.Range("S5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/artfynd/A 1510-2024 artfynd.xlsx", "A 1510-2024")'
$ws.Range("T5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/kartor/A 1510-2024 karta.png", "A 1510-2024")'
$ws.Range("V5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/klagomål/A 1510-2024 FSC-klagomål.docx", "A 1510-2024")'
$ws.Range("W5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/klagomålsmail/A 1510-2024 FSC-klagomål mail.docx", "A 1510-2024")'
$ws.Range("X5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/tillsyn/A 1510-2024 tillsynsbegäran.docx", "A 1510-2024")'
$ws.Range("Y5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/tillsynsmail/A 1510-2024 tillsynsbegäran mail.docx", "A 1510-2024")'
$ws.Range("Z5").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/fåglar/A 1510-2024 prioriterade fågelarter.docx", "A 1510-2024")'
# --- row 6 (source row 4) ---
$ws.Range("A6").Value = 'A 24809-2025'
$ws.Range("B6").Value = 45799
$ws.Range("C6").Value = 46063
$ws.Range("G6").Value = 1.4
$ws.Range("J6").Value = 1
$ws.Range("N6").Value = 0
$ws.Range("R6").Value = 'Spillkråka'
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/artfynd/A 24809-2025 artfynd.xlsx", "A 24809-2025")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/kartor/A 24809-2025 karta.png", "A 24809-2025")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/klagomål/A 24809-2025 FSC-klagomål.docx", "A 24809-2025")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/klagomålsmail/A 24809-2025 FSC-klagomål mail.docx", "A 24809-2025")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/tillsyn/A 24809-2025 tillsynsbegäran.docx", "A 24809-2025")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/tillsynsmail/A 24809-2025 tillsynsbegäran mail.docx", "A 24809-2025")'
$ws.Range("Z6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_1285/fåglar/A 24809-2025 prioriterade fågelarter.docx", "A 24809-2025")'
# --- row 7 (source row 7) ---
$ws.Range("C7").Value = 46063
# --- row 8 (source row 8) ---
$ws.Range("C8").Value = 46063
# --- row 9 (source row 9) ---
$ws.Range("C9").Value = 46063
# --- row 10 (source row 38) ---
$ws.Range("A10").Value = 'A 15533-2024'
$ws.Range("B10").Value = 45401
$ws.Range("C10").Value = 46063
$ws.Range("G10").Value = 7.5
# --- row 11 (source row 51) ---
$ws.Range("A11").Value = 'A 11221-2024'
$ws.Range("B11").Value = 45371.59097222222
$ws.Range("C11").Value = 46063
$ws.Range("G11").Value = 0.3
# --- row 12 (source row 22) ---
$ws.Range("A12").Value = 'A 15467-2025'
$ws.Range("B12").Value = 45747.52465277778
$ws.Range("C12").Value = 46063
$ws.Range("G12").Value = 1.4
# --- row 13 (source row 33) ---
$ws.Range("A13").Value = 'A 12079-2024'
$ws.Range("B13").Value = 45377.542291666665
$ws.Range("C13").Value = 46063
$ws.Range("G13").Value = 6.9
# --- row 14 (source row 25) ---
$ws.Range("A14").Value = 'A 5869-2023'
$ws.Range("B14").Value = 44958
$ws.Range("C14").Value = 46063
$ws.Range("G14").Value = 3
# --- row 15 (source row 37) ---
$ws.Range("A15").Value = 'A 9482-2023'
$ws.Range("B15").Value = 44981
$ws.Range("C15").Value = 46063
$ws.Range("G15").Value = 3.7
# --- row 16 (source row 41) ---
$ws.Range("A16").Value = 'A 25385-2023'
$ws.Range("B16").Value = 45089.337118055555
$ws.Range("C16").Value = 46063
$ws.Range("G16").Value = 0.7
# --- row 17 (source row 52) ---
$ws.Range("A17").Value = 'A 43808-2024'
$ws.Range("B17").Value = 45571
$ws.Range("C17").Value = 46063
$ws.Range("G17").Value = 1.5
# --- row 18 (source row 46) ---
$ws.Range("A18").Value = 'A 12605-2025'
$ws.Range("B18").Value = 45733
$ws.Range("C18").Value = 46063
$ws.Range("F18").Value = 'Sveaskog'
$ws.Range("G18").Value = 2.4
# --- row 19 (source row 34) ---
$ws.Range("A19").Value = 'A 11204-2023'
$ws.Range("B19").Value = 44992.60141203704
$ws.Range("C19").Value = 46063
$ws.Range("G19").Value = 0.7
# --- row 20 (source row 60) ---
$ws.Range("A20").Value = 'A 25508-2023'
$ws.Range("B20").Value = 45089
$ws.Range("C20").Value = 46063
$ws.Range("G20").Value = 0.4
# --- row 21 (source row 27) ---
$ws.Range("A21").Value = 'A 9686-2023'
$ws.Range("B21").Value = 44984
$ws.Range("C21").Value = 46063
$ws.Range("G21").Value = 0.8
# --- row 22 (source row 17) ---
$ws.Range("A22").Value = 'A 41803-2025'
$ws.Range("B22").Value = 45902.63017361111
$ws.Range("C22").Value = 46063
$ws.Range("G22").Value = 1.3
# --- row 23 (source row 35) ---
$ws.Range("A23").Value = 'A 61020-2024'
$ws.Range("B23").Value = 45645.45321759259
$ws.Range("C23").Value = 46063
$ws.Range("G23").Value = 4.9
# --- row 24 (source row 19) ---
$ws.Range("A24").Value = 'A 43943-2025'
$ws.Range("B24").Value = 45915.366273148145
$ws.Range("C24").Value = 46063
$ws.Range("G24").Value = 6.1
# --- row 25 (source row 26) ---
$ws.Range("A25").Value = 'A 7409-2025'
$ws.Range("B25").Value = 45705.36702546296
$ws.Range("C25").Value = 46063
$ws.Range("G25").Value = 1.5
# --- row 26 (source row 23) ---
$ws.Range("A26").Value = 'A 47345-2025'
$ws.Range("B26").Value = 45930.57266203704
$ws.Range("C26").Value = 46063
$ws.Range("G26").Value = 0.7
# --- row 27 (source row 24) ---
$ws.Range("A27").Value = 'A 47356-2025'
$ws.Range("B27").Value = 45930
$ws.Range("C27").Value = 46063
$ws.Range("G27").Value = 0.6
# --- row 28 (source row 45) ---
$ws.Range("A28").Value = 'A 36761-2022'
$ws.Range("B28").Value = 44805
$ws.Range("C28").Value = 46063
$ws.Range("F28").ClearContents() | Out-Null
$ws.Range("G28").Value = 4.8
# --- row 29 (source row 61) ---
$ws.Range("A29").Value = 'A 42050-2023'
$ws.Range("B29").Value = 45177
$ws.Range("C29").Value = 46063
$ws.Range("G29").Value = 10.3
# --- row 30 (source row 32) ---
$ws.Range("A30").Value = 'A 25388-2023'
$ws.Range("B30").Value = 45089.342199074075
$ws.Range("C30").Value = 46063
$ws.Range("G30").Value = 0.5
# --- row 31 (source row 30) ---
$ws.Range("A31").Value = 'A 35996-2025'
$ws.Range("C31").Value = 46063
$ws.Range("G31").Value = 0.8
# --- row 32 (source row 48) ---
$ws.Range("A32").Value = 'A 6545-2023'
$ws.Range("B32").Value = 44960
$ws.Range("C32").Value = 46063
$ws.Range("G32").Value = 4.1
# --- row 33 (source row 31) ---
$ws.Range("A33").Value = 'A 35997-2025'
$ws.Range("B33").Value = 45863
$ws.Range("C33").Value = 46063
$ws.Range("G33").Value = 0.9
# --- row 34 (source row 12) ---
$ws.Range("A34").Value = 'A 13437-2021'
$ws.Range("B34").Value = 44273
$ws.Range("C34").Value = 46063
$ws.Range("G34").Value = 1.3
# --- row 35 (source row 10) ---
$ws.Range("A35").Value = 'A 61035-2024'
$ws.Range("B35").Value = 45645.46634259259
$ws.Range("C35").Value = 46063
$ws.Range("G35").Value = 7.6
# --- row 36 (source row 50) ---
$ws.Range("A36").Value = 'A 10940-2024'
$ws.Range("B36").Value = 45370.45334490741
$ws.Range("C36").Value = 46063
$ws.Range("G36").Value = 0.5
# --- row 37 (source row 49) ---
$ws.Range("A37").Value = 'A 11732-2025'
$ws.Range("B37").Value = 45727
$ws.Range("C37").Value = 46063
$ws.Range("G37").Value = 4.5
# --- row 38 (source row 28) ---
$ws.Range("A38").Value = 'A 7016-2025'
$ws.Range("B38").Value = 45701
$ws.Range("C38").Value = 46063
$ws.Range("F38").Value = 'Sveaskog'
$ws.Range("G38").Value = 2.3
# --- row 39 (source row 36) ---
$ws.Range("A39").Value = 'A 2691-2026'
$ws.Range("B39").Value = 46037.66427083333
$ws.Range("C39").Value = 46063
$ws.Range("G39").Value = 0.8
# --- row 40 (source row 40) ---
$ws.Range("C40").Value = 46063
# --- row 41 (source row 20) ---
$ws.Range("A41").Value = 'A 20160-2024'
$ws.Range("B41").Value = 45434
$ws.Range("C41").Value = 46063
$ws.Range("G41").Value = 4.1
# --- row 42 (source row 62) ---
$ws.Range("A42").Value = 'A 14488-2024'
$ws.Range("B42").Value = 45394
$ws.Range("C42").Value = 46063
$ws.Range("F42").Value = 'Sveaskog'
$ws.Range("G42").Value = 1.6
# --- row 43 (source row 57) ---
$ws.Range("A43").Value = 'A 16670-2025'
$ws.Range("B43").Value = 45754.43791666667
$ws.Range("C43").Value = 46063
$ws.Range("G43").Value = 1.5
# --- row 44 (source row 47) ---
$ws.Range("A44").Value = 'A 34468-2023'
$ws.Range("B44").Value = 45139
$ws.Range("C44").Value = 46063
$ws.Range("G44").Value = 5.5
# --- row 45 (source row 43) ---
$ws.Range("A45").Value = 'A 47870-2025'
$ws.Range("B45").Value = 45932
$ws.Range("C45").Value = 46063
$ws.Range("G45").Value = 6.8
# --- row 46 (source row 11) ---
$ws.Range("A46").Value = 'A 10773-2025'
$ws.Range("B46").Value = 45722.47409722222
$ws.Range("C46").Value = 46063
$ws.Range("F46").ClearContents() | Out-Null
$ws.Range("G46").Value = 1.2
# --- row 47 (source row 39) ---
$ws.Range("A47").Value = 'A 62496-2023'
$ws.Range("B47").Value = 45268
$ws.Range("C47").Value = 46063
$ws.Range("G47").Value = 2.6
# --- row 48 (source row 42) ---
$ws.Range("A48").Value = 'A 34170-2023'
$ws.Range("B48").Value = 45138
$ws.Range("C48").Value = 46063
$ws.Range("G48").Value = 5.8
# --- row 49 (source row 29) ---
$ws.Range("A49").Value = 'A 24115-2022'
$ws.Range("B49").Value = 44725
$ws.Range("C49").Value = 46063
$ws.Range("G49").Value = 4.4
# --- row 50 (source row 53) ---
$ws.Range("A50").Value = 'A 9484-2023'
$ws.Range("B50").Value = 44981
$ws.Range("C50").Value = 46063
$ws.Range("G50").Value = 1.4
# --- row 51 (source row 54) ---
$ws.Range("A51").Value = 'A 9487-2023'
$ws.Range("B51").Value = 44981
$ws.Range("C51").Value = 46063
$ws.Range("G51").Value = 1.9
# --- row 52 (source row 14) ---
$ws.Range("A52").Value = 'A 16733-2023'
$ws.Range("B52").Value = 45030
$ws.Range("C52").Value = 46063
$ws.Range("G52").Value = 2.6
# --- row 53 (source row 55) ---
$ws.Range("A53").Value = 'A 44496-2023'
$ws.Range("B53").Value = 45189
$ws.Range("C53").Value = 46063
$ws.Range("G53").Value = 0.5
# --- row 54 (source row 58) ---
$ws.Range("A54").Value = 'A 34466-2023'
$ws.Range("B54").Value = 45139
$ws.Range("C54").Value = 46063
$ws.Range("G54").Value = 4.3
# --- row 55 (source row 56) ---
$ws.Range("A55").Value = 'A 33865-2024'
$ws.Range("B55").Value = 45520
$ws.Range("C55").Value = 46063
$ws.Range("G55").Value = 1.3
# --- row 56 (source row 18) ---
$ws.Range("A56").Value = 'A 5812-2022'
$ws.Range("B56").Value = 44596
$ws.Range("C56").Value = 46063
$ws.Range("G56").Value = 0.5
# --- row 57 (source row 59) ---
$ws.Range("A57").Value = 'A 38173-2023'
$ws.Range("B57").Value = 45161.44363425926
$ws.Range("C57").Value = 46063
$ws.Range("G57").Value = 2.8
# --- row 58 (source row 44) ---
$ws.Range("A58").Value = 'A 15357-2023'
$ws.Range("B58").Value = 45019
$ws.Range("C58").Value = 46063
$ws.Range("G58").Value = 1
# --- row 59 (source row 21) ---
$ws.Range("A59").Value = 'A 1574-2024'
$ws.Range("B59").Value = 45306
$ws.Range("C59").Value = 46063
$ws.Range("G59").Value = 3.6
# --- row 60 (source row 63) ---
$ws.Range("A60").Value = 'A 49137-2024'
$ws.Range("B60").Value = 45594
$ws.Range("C60").Value = 46063
$ws.Range("G60").Value = 2.4
# --- row 61 (source row 64) ---
$ws.Range("A61").Value = 'A 52652-2023'
$ws.Range("B61").Value = 45225
$ws.Range("C61").Value = 46063
$ws.Range("G61").Value = 2.9
# --- row 62 (source row 66) ---
$ws.Range("A62").Value = 'A 32984-2024'
$ws.Range("B62").Value = 45517
$ws.Range("C62").Value = 46063
$ws.Range("F62").ClearContents() | Out-Null
$ws.Range("G62").Value = 2.9
# --- row 63 (source row 65) ---
$ws.Range("A63").Value = 'A 18713-2024'
$ws.Range("B63").Value = 45426
$ws.Range("C63").Value = 46063
$ws.Range("G63").Value = 4.1
# --- row 64 (source row 67) ---
$ws.Range("A64").Value = 'A 23767-2025'
$ws.Range("B64").Value = 45793.47238425926
$ws.Range("C64").Value = 46063
$ws.Range("G64").Value = 3.2
# --- row 65 (source row 68) ---
$ws.Range("A65").Value = 'A 23773-2025'
$ws.Range("B65").Value = 45793.48923611111
$ws.Range("C65").Value = 46063
$ws.Range("G65").Value = 4.9
# --- row 66 (source row 13) ---
$ws.Range("A66").Value = 'A 16792-2024'
$ws.Range("B66").Value = 45411
$ws.Range("C66").Value = 46063
$ws.Range("G66").Value = 0.9
# --- row 67 (source row 16) ---
$ws.Range("A67").Value = 'A 28088-2025'
$ws.Range("B67").Value = 45817
$ws.Range("C67").Value = 46063
$ws.Range("G67").Value = 1
# --- row 68 (source row 15) ---
$ws.Range("A68").Value = 'A 28093-2025'
$ws.Range("B68").Value = 45817
$ws.Range("C68").Value = 46063
$ws.Range("G68").Value = 10.2
